# Fix the PER (Player Efficiency Rating) bug: the team order/labels in
# column B were shuffled relative to the shared-string table, and the C
# column held the wrong (weight-like) numbers instead of the real PER
# values. Re-write both columns, row by row, with the corrected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; Team = "POR"; Val = 13.61538461538461 },
    @{ Row = 3; Team = "NJN"; Val = 13.91666666666667 },
    @{ Row = 4; Team = "CLE"; Val = 15.125 },
    @{ Row = 5; Team = "DAL"; Val = 11.96428571428571 },
    @{ Row = 6; Team = "MIA"; Val = 13.475 },
    @{ Row = 7; Team = "SEA"; Val = 13.43846153846154 },
    @{ Row = 8; Team = "ATL"; Val = 15.15714285714285 },
    @{ Row = 9; Team = "MIL"; Val = 13.59230769230769 },
    @{ Row = 10; Team = "LAC"; Val = 12.61538461538461 },
    @{ Row = 11; Team = "DET"; Val = 12.50769230769231 },
    @{ Row = 12; Team = "SAS"; Val = 14.87857142857143 },
    @{ Row = 13; Team = "ORL"; Val = 13.01428571428571 },
    @{ Row = 14; Team = "UTA"; Val = 13.26666666666667 },
    @{ Row = 15; Team = "HOU"; Val = 11.94615384615385 },
    @{ Row = 16; Team = "DEN"; Val = 12.52307692307692 },
    @{ Row = 17; Team = "LAL"; Val = 12.88461538461539 },
    @{ Row = 18; Team = "GSW"; Val = 12.575 },
    @{ Row = 19; Team = "IND"; Val = 13.23076923076923 },
    @{ Row = 20; Team = "CHI"; Val = 13.15 },
    @{ Row = 21; Team = "PHI"; Val = 11.43076923076923 },
    @{ Row = 22; Team = "CHH"; Val = 12.8076923076923 },
    @{ Row = 23; Team = "BOS"; Val = 13.23571428571428 },
    @{ Row = 24; Team = "WSB"; Val = 11.91333333333333 },
    @{ Row = 25; Team = "SAC"; Val = 12.30833333333333 },
    @{ Row = 26; Team = "PHO"; Val = 14.71538461538462 },
    @{ Row = 27; Team = "NYK"; Val = 10.44166666666667 },
    @{ Row = 28; Team = "MIN"; Val = 13.39166666666667 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.Team
    $ws.Cells.Item($r, 3).Value = $item.Val
}
